# Test_Cases.xlsx - "Creation of Ticket" test rows: fill in the Actual Output /
# Pass columns now that testing has been performed.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6: Creation of Ticket by User (all fields entered) -> passed
$ws.Range("E6").Value = "User able to submit ticket"
$ws.Range("F6").Value = "Yes"

# Row 7: Creation of Ticket by Tech (all fields entered) -> passed
$ws.Range("E7").Value = "Tech able to submit ticket"
$ws.Range("F7").Value = "Yes"

# Row 8: Invalid/Missing inputs on ticket creation -> bug found, did not pass
# (a ticket was created even though required data was missing)
$ws.Range("E8").Value = "Tickets were created"
$ws.Range("F8").Value = "No"

# Leave the selection where the author ended up while entering this data
$ws.Range("E12").Select()
